$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)   # "Overview"
$zhcn     = $wb.Worksheets.Item(2)   # "zh-cn"
$dede     = $wb.Worksheets.Item(3)   # "de-de"

$newStatus = "Handed back: in sync with en-US"

# ----- Overview sheet: Status columns for zh-cn (E) and de-de (F) -----
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# ----- zh-cn sheet -----
$zhcn.Range("C2").Value = $newStatus                      # Status
$zhcn.Range("K2").Value = "2016-09-06 15:43:42"            # Latest Handback DateTime
$zhcn.Range("P2").ClearContents()                          # Error Detail no longer applies

# ----- de-de sheet -----
$dede.Range("C2").Value = $newStatus                       # Status
$dede.Range("K2").Value = "2016-09-06 15:43:50"             # Latest Handback DateTime
$dede.Range("P2").ClearContents()                           # Error Detail no longer applies

# ----- Column width adjustments (status columns grew wider, error-detail columns shrank) -----
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527    # E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527    # F (de-de status)

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527        # C (Status)
$zhcn.Columns.Item(16).ColumnWidth = 13.7470528738839        # P (Error Detail)

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527        # C (Status)
$dede.Columns.Item(16).ColumnWidth = 13.7470528738839        # P (Error Detail)
